$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '67.210.41'
$ws.Cells.Item(2, 5).Value = '  -0.79%  '
$ws.Cells.Item(3, 4).Value = '3.512.81'
$ws.Cells.Item(3, 5).Value = '  +0.18%  '
$ws.Cells.Item(4, 5).Value = '  -0.20%  '
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '610.65'
$ws.Cells.Item(5, 5).Value = '  +0.47%  '
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '148.14'
$ws.Cells.Item(6, 5).Value = '  -1.61%  '
$ws.Cells.Item(7, 4).Value = '3.512.59'
$ws.Cells.Item(7, 5).Value = '  +0.30%  '
$ws.Cells.Item(8, 5).Value = '  -0.04%  '
$ws.Cells.Item(9, 5).Value = '  -1.46%  '
$ws.Cells.Item(10, 5).Value = '  -1.19%  '
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = '8.07'
$ws.Cells.Item(11, 5).Value = '  +6.48%  '
$ws.Cells.Item(12, 5).Value = '  -1.68%  '
$ws.Cells.Item(13, 5).Value = '  +0.97%  '
$ws.Cells.Item(14, 4).Value = '4.095.71'
$ws.Cells.Item(14, 5).Value = '  -0.21%  '
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = '31.56'
$ws.Cells.Item(15, 5).Value = '  -1.63%  '
$ws.Cells.Item(16, 4).Value = '3.509.82'
$ws.Cells.Item(16, 5).Value = '  +0.29%  '
$ws.Cells.Item(17, 4).Value = '67.236.00'
$ws.Cells.Item(17, 5).Value = '  -1.22%  '
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = '0.117'
$ws.Cells.Item(18, 5).Value = '  -0.19%  '
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '10.90'
$ws.Cells.Item(19, 5).Value = '  +9.62%  '
$ws.Cells.Item(20, 5).Value = '  -2.31%  '
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = '15.43'
$ws.Cells.Item(21, 5).Value = '  +0.11%  '
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '436.53'
$ws.Cells.Item(22, 5).Value = '  -2.80%  '
$ws.Cells.Item(23, 5).Value = '  -2.63%  '
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = '80.11'
$ws.Cells.Item(24, 5).Value = '  +1.77%  '
$ws.Cells.Item(25, 4).Value = '3.651.34'
$ws.Cells.Item(25, 5).Value = '  +0.03%  '
$ws.Cells.Item(26, 5).Value = '  +0.02%  '
$ws.Cells.Item(27, 5).Value = '  -3.67%  '
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = '9.87'
$ws.Cells.Item(28, 5).Value = '  -0.69%  '
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = '8.29'
$ws.Cells.Item(29, 5).Value = '  -4.77%  '
$ws.Cells.Item(30, 5).Value = '  +0.34%  '
$ws.Cells.Item(31, 5).Value = '  -3.95%  '
$ws.Cells.Item(32, 5).Value = '  +0.02%  '
$ws.Cells.Item(33, 5).Value = '  -2.00%  '
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = '25.60'
$ws.Cells.Item(34, 5).Value = '  -0.04%  '
$ws.Cells.Item(35, 5).Value = '  -1.49%  '
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = '5.96'
$ws.Cells.Item(36, 5).Value = '  -3.79%  '
$ws.Cells.Item(37, 5).Value = '  +0.09%  '
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '0.999'
$ws.Cells.Item(39, 5).Value = '  -0.21%  '
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '175.59'
$ws.Cells.Item(40, 5).Value = '  -1.89%  '
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '0.0905'
$ws.Cells.Item(41, 5).Value = '  -0.24%  '
$ws.Cells.Item(42, 5).Value = '  -0.14%  '
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = '2.05'
$ws.Cells.Item(43, 5).Value = '  -9.86%  '
$ws.Cells.Item(44, 5).Value = '  -0.12%  '
$ws.Cells.Item(45, 5).Value = '  -1.25%  '
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '28.70'
$ws.Cells.Item(46, 5).Value = '  -8.40%  '
$ws.Cells.Item(47, 5).Value = '  -5.54%  '
$ws.Cells.Item(48, 5).Value = '  -2.02%  '
$ws.Cells.Item(49, 5).Value = '  -2.27%  '
$ws.Cells.Item(50, 5).Value = '  -1.20%  '
$ws.Cells.Item(51, 5).Value = '  -1.92%  '
